# population_share_children.xlsx — strip the two-space indent that was used
# to visually "nest" district names under their province in column B, now
# that column A/B carry province/district as plain values, and left-align
# the whole used range (was center/right for the header + number columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Clean up district names in column B (drop the leading "  ") ------
$districtRows = @(4,5,6,8,9,10,11,12,13,14,15,17,18,19,20,21,22,23,25,26,27,28,29,31,32,33,34,35,36,37)
foreach ($r in $districtRows) {
    $cell = $ws.Cells.Item($r, 2)
    $cell.Value = $cell.Value2.Trim()
}

# --- 2. Left-align everything (was: general/center header, right numbers) -
$ws.Range("A1:E37").HorizontalAlignment = -4131

# --- 3. Give columns C:E their (now left-aligned) default width so new ---
#        rows added below the data keep the same look.
$ws.Columns("C:E").ColumnWidth = 8.88671875

# --- 4. Move the active selection the way the author left it -------------
$ws.Range("H6").Select()
